# Apply the edit described by the diff:
#  - Insert a new order row (mã dịch vụ 618, khách hàng "đường thị út",
#    dịch vụ "Cắt mí") above the existing "619 / mai hồng nương / Thu cánh mũi"
#    row on sheet "Đơn sale chính".
#  - Recompute the "Tổng" (totals) row to include both orders.
#  - Propagate the new "Chiết khấu sale chính tại SÓC TRĂNG" total (N4 = 1380000)
#    to the corresponding rows on the "Lương" sheet.

$wb = $excel.ActiveWorkbook

$wsOrders = $wb.Worksheets.Item("Đơn sale chính")
$wsSalary = $wb.Worksheets.Item("Lương")

# --- Sheet "Đơn sale chính" -------------------------------------------------

# Insert a new blank row at row 2; existing row 2 (619 ...) and row 3 (Tổng)
# shift down to rows 3 and 4 respectively.
$wsOrders.Rows.Item(2).Insert()

# Fill the new row 2 with the new order's data.
$wsOrders.Range("A2").Value2 = "HD-LUXURY"
$wsOrders.Range("B2").Value2 = 618
# Force the date column to stay a plain text string (matches the rest of the
# sheet, which stores "Ngày thực hiện" as inline text rather than a date) -
# temporarily apply a text format so Excel doesn't auto-convert the literal
# into a date serial number, then drop back to the default "Normal" style so
# no stray cell formatting is left behind.
$wsOrders.Range("C2").NumberFormat = "@"
$wsOrders.Range("C2").Value2 = "08-02-2024"
$wsOrders.Range("C2").Style = "Normal"
$wsOrders.Range("D2").Value2 = "SÓC TRĂNG"
$wsOrders.Range("E2").Value2 = "đường thị út"
$wsOrders.Range("F2").Value2 = "Cá nhân"
$wsOrders.Range("G2").Value2 = "Cắt mí"
$wsOrders.Range("H2").Value2 = 6000000
$wsOrders.Range("I2").Value2 = $null
$wsOrders.Range("J2").Value2 = $null
$wsOrders.Range("K2").Value2 = 6000000
$wsOrders.Range("L2").Value2 = 6000000
$wsOrders.Range("M2").Value2 = 0.13
$wsOrders.Range("N2").Value2 = 780000

# Row 3 now holds the original order (619 / mai hồng nương / Thu cánh mũi),
# whose M3 discount ratio was previously stored on the totals row; restore it
# on the data row itself.
$wsOrders.Range("M3").Value2 = 0.1

# Keep the "Sale phụ" / "Upsale" columns blank (numeric, empty) on the data
# row, same as before the insert shifted row 2 down to row 3.
$wsOrders.Range("I3").Value2 = $null
$wsOrders.Range("J3").Value2 = $null

# Row 4 is now the totals ("Tổng") row; recompute its aggregated values.
$wsOrders.Range("A4").Value2 = "Tổng"
$wsOrders.Range("B4").Value2 = 2
$wsOrders.Range("H4").Value2 = 14000000
$wsOrders.Range("J4").Value2 = 0
$wsOrders.Range("K4").Value2 = 14000000
$wsOrders.Range("L4").Value2 = 12000000
$wsOrders.Range("M4").Value2 = 0
$wsOrders.Range("N4").Value2 = 1380000

# --- Sheet "Lương" -----------------------------------------------------------

# The "Chiết khấu sale chính tại SÓC TRĂNG" total, and the totals derived from
# it, grow from 600000 to 1380000 (matching N4 above).
$wsSalary.Range("B25").Value2 = 1380000
$wsSalary.Range("B34").Value2 = 1380000
$wsSalary.Range("B35").Value2 = 1380000
